# DesignerConfigs/Datas/test/global.xlsx
#
# luban-server no longer supports the orientation=landscape|l|portrait|p
# meta values, so the "orientation=portrait" comment cell is rewritten to
# "orientation=c" (and loses its mixed-run rich formatting in the process,
# since the new value is written as a single plain run).
#
# The active selection on the TbGlobalConfig sheet is also moved up one
# row, from B11 to B10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TbGlobalConfig")

# B1 holds the "orientation=portrait" note (stored as a shared string with
# two formatting runs). Replace its text outright with the new value.
$ws.Range("B1").Value = "orientation=c"

# Move the saved cell selection from B11 to B10.
$ws.Range("B10").Select()
